$wb = $excel.ActiveWorkbook

# --- Reorganise sheets: insert a new "2022-Q1" sheet between "2021-Q4" and
# "总计", and make sure sheetId order matches (1, 2, 3). We do this by
# temporarily removing "总计", adding "2022-Q1" right after "2021-Q4", then
# re-adding "总计" (so it gets the next free sheetId) and moving it back to
# the end.

$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Delete()

$wsQ4 = $wb.Worksheets.Item(1)
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# --- Populate the new "2022-Q1" sheet with the per-fund holdings table.

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $wsQ1.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Style = "Bold Centered"
}

$rows = @(
    @("398001", "中海优质成长混合", "14.42", "90.86", "3.74", "0.5393", 8),
    @("002430", "中银丰利灵活配置混合A", "8.76", "28.78", "1.92", "0.1682", 9),
    @("002616", "中银益利灵活配置混合A", "5.90", "29.85", "2.22", "0.1310", 5),
    @("003850", "中银锦利灵活配置混合A", "6.37", "28.68", "1.97", "0.1255", 8),
    @("001370", "中银新趋势灵活配置混合", "2.63", "34.29", "2.33", "0.0613", 7),
    @("002431", "中银丰利灵活配置混合C", "2.24", "28.78", "1.92", "0.0430", 9),
    @("003851", "中银锦利灵活配置混合C", "1.25", "28.68", "1.97", "0.0246", 8),
    @("002617", "中银益利灵活配置混合C", "0.58", "29.85", "2.22", "0.0129", 5),
    @("011677", "中银睿丰回报混合型证券投资基金A", "0.73", "20.29", "1.55", "0.0113", 9),
    @("001252", "中海进取收益灵活配置混合", "0.23", "92.65", "3.92", "0.0090", 8),
    @("011678", "中银睿丰回报混合型证券投资基金C", "0.00", "20.29", "1.55", 0, 9)
)

$r = 2
foreach ($row in $rows) {
    $wsQ1.Cells.Item($r, 1).Value = $r - 2
    $wsQ1.Cells.Item($r, 1).Style = "Bold Centered"
    $wsQ1.Cells.Item($r, 2).Value = $row[0]
    $wsQ1.Cells.Item($r, 3).Value = $row[1]
    $wsQ1.Cells.Item($r, 4).Value = $row[2]
    $wsQ1.Cells.Item($r, 5).Value = $row[3]
    $wsQ1.Cells.Item($r, 6).Value = $row[4]
    $wsQ1.Cells.Item($r, 7).Value = $row[5]
    $wsQ1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# --- Insert the new "2022-Q1" summary row into "总计", pushing the
# existing "2021-Q4" row down.

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 1).Style = "Bold Centered"
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 11
$wsTotal.Cells.Item(2, 4).Value = 1.13

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 1).Style = "Bold Centered"
